$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-14 per repull/recalculation of data
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = -1
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = -5
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = 7
